# Update "想去人数" (interested-people count) figures that changed when the
# site was regenerated, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1128
$ws1.Range("F8").Value  = 11508
$ws1.Range("F12").Value = 20
$ws1.Range("F14").Value = 2531
$ws1.Range("F16").Value = 129
$ws1.Range("F18").Value = 2210
$ws1.Range("F20").Value = 510
$ws1.Range("F21").Value = 11294
$ws1.Range("F22").Value = 11204

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1128
$ws4.Range("F8").Value  = 11508
$ws4.Range("F12").Value = 20
$ws4.Range("F14").Value = 2531
$ws4.Range("F17").Value = 129
$ws4.Range("F19").Value = 2210
$ws4.Range("F21").Value = 510
$ws4.Range("F22").Value = 11294
$ws4.Range("F23").Value = 11204
